{"js": "async (context) => {\n  // The edit removes the stray \" 's\" typo in \"it's heightens\" -> \"it heightens\"\n  // within the \"My Reflection\" paragraph of the journal.\n  const results = context.document.body.search(\"it\\u2019s heightens the complexity a bit\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"it heightens the complexity a bit\", \"Replace\");\n  }\n  await context.sync();\n};\n", "ps1": "# Fix the stray apostrophe-s typo in the \"My Reflection\" paragraph:\n# \"That's a good thing, it's heightens the complexity a bit.\" ->\n# \"That's a good thing, it heightens the complexity a bit.\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"it\" + [char]0x2019 + \"s heightens\"\n$find.Replacement.Text = \"it heightens\"\n\n# 0 = wdFindContinue (no wrap behavior needed here, searching whole story),\n# 2 = wdReplaceAll\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nWrite-Output \"done\"\n"}
